$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values (column B)
$ws.Range("B2").Value = 102
$ws.Range("B3").Value = 97
$ws.Range("B4").Value = 94
$ws.Range("B5").Value = 83
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 78
$ws.Range("B8").Value = 78
$ws.Range("B9").Value = 77
$ws.Range("B10").Value = 63

# Swap names in rows 7 and 8 (column A)
$ws.Range("A7").Value = "BLANCO LOZANO ANDREA MILAGROS"
$ws.Range("A8").Value = "BURGA MEDINA SHIRLEY ROCIO"
